$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G4").Value = 2.75
$ws.Range("I4").Value = 2.45
$ws.Range("AD5").Value = 9.5
$ws.Range("AG5").Value = 26
$ws.Range("AM5").Value = 67
$ws.Range("AN5").Value = 34
$ws.Range("M7").Value = 1.11
$ws.Range("N7").Value = 6.5
$ws.Range("AR7").Value = 1.98
$ws.Range("AS7").Value = 1.88
$ws.Range("Q8").Value = 2.1
$ws.Range("R8").Value = 1.73
$ws.Range("U8").Value = 3.75
$ws.Range("V8").Value = 1.29
$ws.Range("G10").Value = 1.6
$ws.Range("H10").Value = 3.7
$ws.Range("I10").Value = 5.5
$ws.Range("J10").Value = 2.25
$ws.Range("K10").Value = 2.1
$ws.Range("Y10").Value = 2.05
$ws.Range("Z10").Value = 1.7
$ws.Range("AB10").Value = 7
$ws.Range("AD10").Value = 12
$ws.Range("AF10").Value = 29
$ws.Range("AH10").Value = 7
$ws.Range("AI10").Value = 19
$ws.Range("AK10").Value = 451
$ws.Range("AL10").Value = 13
$ws.Range("AP10").Value = 41
$ws.Range("L17").Value = 1.91
$ws.Range("O17").Value = 1.17
$ws.Range("P17").Value = 4.5
$ws.Range("W17").Value = 1.29
$ws.Range("X17").Value = 3.5
$ws.Range("AL17").Value = 9
$ws.Range("AN17").Value = 8.5
$ws.Range("G20").Value = 3.2
$ws.Range("H20").Value = 3.4
$ws.Range("I20").Value = 2.2
$ws.Range("J20").Value = 3.5
$ws.Range("K20").Value = 2.2
$ws.Range("L20").Value = 2.75
$ws.Range("M20").Value = 1.02
$ws.Range("N20").Value = 11
$ws.Range("W20").Value = 1.36
$ws.Range("X20").Value = 3
$ws.Range("Y20").Value = 1.67
$ws.Range("Z20").Value = 2.1
$ws.Range("AC20").Value = 12
$ws.Range("AE20").Value = 23
$ws.Range("AF20").Value = 29
$ws.Range("AG20").Value = 11
$ws.Range("AL20").Value = 9
$ws.Range("AM20").Value = 11
$ws.Range("AN20").Value = 9.5
$ws.Range("AO20").Value = 21
$ws.Range("J23").Value = 4.2
$ws.Range("L23").Value = 2.37
$ws.Range("W23").Value = 1.36
$ws.Range("X23").Value = 3
$ws.Range("AA23").Value = 12.5
$ws.Range("AC23").Value = 14
$ws.Range("AE23").Value = 37
$ws.Range("AF23").Value = 37
$ws.Range("AK23").Value = 400
$ws.Range("AL23").Value = 8.25
$ws.Range("AM23").Value = 10.25
$ws.Range("AQ23").Value = 23
